$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetOld = "System, dnasr281@gmail.com"
$targetNew = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $targetOld) {
        $cell.Value2 = $targetNew
    }
}
